$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = '="58.542.11"'
$ws.Range("E2").Formula = '="  -2.69%  "'
$ws.Range("D3").Formula = '="2.276.13"'
$ws.Range("E3").Formula = '="  -6.04%  "'
$ws.Range("D4").Formula = '="0.999"'
$ws.Range("E4").Formula = '="  -0.09%  "'
$ws.Range("D5").Formula = '="543.75"'
$ws.Range("E5").Formula = '="  -1.70%  "'
$ws.Range("E6").Formula = '="  -5.08%  "'
$ws.Range("E7").Formula = '="  -0.03%  "'
$ws.Range("E8").Formula = '="  -3.62%  "'
$ws.Range("E9").Formula = '="  -4.09%  "'
$ws.Range("D10").Formula = '="5.50"'
$ws.Range("E10").Formula = '="  -3.38%  "'
$ws.Range("D11").Formula = '="0.149"'
$ws.Range("E11").Formula = '="  +0.54%  "'
$ws.Range("E12").Formula = '="  -5.73%  "'
$ws.Range("D13").Formula = '="23.65"'
$ws.Range("E13").Formula = '="  -5.66%  "'
$ws.Range("D14").Formula = '="2.681.01"'
$ws.Range("E14").Formula = '="  -6.09%  "'
$ws.Range("D15").Formula = '="58.502.08"'
$ws.Range("E15").Formula = '="  -2.62%  "'
$ws.Range("E16").Formula = '="  -3.75%  "'
$ws.Range("D17").Formula = '="2.276.30"'
$ws.Range("E17").Formula = '="  -6.66%  "'
$ws.Range("E18").Formula = '="  -5.82%  "'
$ws.Range("E19").Formula = '="  -4.66%  "'
$ws.Range("D20").Formula = '="314.08"'
$ws.Range("E20").Formula = '="  -4.54%  "'
$ws.Range("E21").Formula = '="  -4.60%  "'
$ws.Range("E22").Formula = '="  +0.08%  "'
$ws.Range("D23").Formula = '="62.86"'
$ws.Range("E23").Formula = '="  -4.32%  "'
$ws.Range("E24").Formula = '="  -4.80%  "'
$ws.Range("D25").Formula = '="1.00"'
$ws.Range("E25").Formula = '="  -0.02%  "'
$ws.Range("E26").Formula = '="  -7.14%  "'
$ws.Range("D27").Formula = '="1.31"'
$ws.Range("E27").Formula = '="  -6.22%  "'
$ws.Range("E28").Formula = '="  -1.66%  "'
$ws.Range("D29").Formula = '="170.73"'
$ws.Range("E29").Formula = '="  +0.37%  "'
$ws.Range("D30").Formula = '="0.0₃0720"'
$ws.Range("E30").Formula = '="  -7.08%  "'
$ws.Range("E31").Formula = '="  -0.27%  "'
$ws.Range("E32").Formula = '="  -5.86%  "'
$ws.Range("E33").Formula = '="  -5.93%  "'
$ws.Range("E35").Formula = '="  -4.64%  "'
$ws.Range("E36").Formula = '="  +0.03%  "'
$ws.Range("E37").Formula = '="  -5.87%  "'
$ws.Range("E38").Formula = '="  -6.58%  "'
$ws.Range("D39").Formula = '="38.11"'
$ws.Range("E39").Formula = '="  -1.89%  "'
$ws.Range("E40").Formula = '="  -5.76%  "'
$ws.Range("D41").Formula = '="300.48"'
$ws.Range("E41").Formula = '="  -10.06%  "'
$ws.Range("D42").Formula = '="140.59"'
$ws.Range("E42").Formula = '="  -3.93%  "'
$ws.Range("E43").Formula = '="  -5.77%  "'
$ws.Range("D44").Formula = '="0.0946"'
$ws.Range("E44").Formula = '="  -2.21%  "'
$ws.Range("E45").Formula = '="  -3.73%  "'
$ws.Range("E46").Formula = '="  -5.06%  "'
$ws.Range("D47").Formula = '="18.32"'
$ws.Range("E47").Formula = '="  -9.07%  "'
$ws.Range("E48").Formula = '="  -4.52%  "'
$ws.Range("E49").Formula = '="  -0.33%  "'
$ws.Range("D50").Formula = '="16.48"'
$ws.Range("E50").Formula = '="  -6.80%  "'
$ws.Range("D51").Formula = '="4.63"'
$ws.Range("E51").Formula = '="  -0.64%  "'

$rng = $ws.Range("D2:E51")
$rng.Copy()
$rng.PasteSpecial(-4163)
$excel.CutCopyMode = $false
